# Edit script: adds a "metadata" sheet after "data" and refreshes the
# "time_taken" timestamps (F2:F65) on the "data" sheet.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Refresh F2:F65 ("time_taken") timestamps on the "data" sheet.
# ---------------------------------------------------------------------------
$newTimestamps = @(
    "2021-10-05 14:35:09.562781",
    "2021-10-05 14:35:09.562788",
    "2021-10-05 14:35:09.562791",
    "2021-10-05 14:35:09.562794",
    "2021-10-05 14:35:09.562797",
    "2021-10-05 14:35:09.562800",
    "2021-10-05 14:35:09.562803",
    "2021-10-05 14:35:09.562805",
    "2021-10-05 14:35:09.562808",
    "2021-10-05 14:35:09.562810",
    "2021-10-05 14:35:09.562813",
    "2021-10-05 14:35:09.562815",
    "2021-10-05 14:35:09.562818",
    "2021-10-05 14:35:09.562820",
    "2021-10-05 14:35:09.562823",
    "2021-10-05 14:35:09.562825",
    "2021-10-05 14:35:09.562828",
    "2021-10-05 14:35:09.562831",
    "2021-10-05 14:35:09.562833",
    "2021-10-05 14:35:09.562836",
    "2021-10-05 14:35:09.562838",
    "2021-10-05 14:35:09.562841",
    "2021-10-05 14:35:09.562843",
    "2021-10-05 14:35:09.562846",
    "2021-10-05 14:35:09.562848",
    "2021-10-05 14:35:09.562851",
    "2021-10-05 14:35:09.562854",
    "2021-10-05 14:35:09.562856",
    "2021-10-05 14:35:09.562859",
    "2021-10-05 14:35:09.562861",
    "2021-10-05 14:35:09.562863",
    "2021-10-05 14:35:09.562866",
    "2021-10-05 14:35:09.562869",
    "2021-10-05 14:35:09.562872",
    "2021-10-05 14:35:09.562874",
    "2021-10-05 14:35:09.562877",
    "2021-10-05 14:35:09.562879",
    "2021-10-05 14:35:09.562882",
    "2021-10-05 14:35:09.562884",
    "2021-10-05 14:35:09.562887",
    "2021-10-05 14:35:09.562890",
    "2021-10-05 14:35:09.562892",
    "2021-10-05 14:35:09.562895",
    "2021-10-05 14:35:09.562897",
    "2021-10-05 14:35:09.562900",
    "2021-10-05 14:35:09.562902",
    "2021-10-05 14:35:09.562905",
    "2021-10-05 14:35:09.562907",
    "2021-10-05 14:35:09.562910",
    "2021-10-05 14:35:09.562912",
    "2021-10-05 14:35:09.562915",
    "2021-10-05 14:35:09.562917",
    "2021-10-05 14:35:09.562920",
    "2021-10-05 14:35:09.562923",
    "2021-10-05 14:35:09.562925",
    "2021-10-05 14:35:09.562928",
    "2021-10-05 14:35:09.562930",
    "2021-10-05 14:35:09.562933",
    "2021-10-05 14:35:09.562935",
    "2021-10-05 14:35:09.562938",
    "2021-10-05 14:35:09.562940",
    "2021-10-05 14:35:09.562943",
    "2021-10-05 14:35:09.562945",
    "2021-10-05 14:35:09.562948"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $dataSheet.Cells.Item($i + 2, 6).Value = $newTimestamps[$i]
}

# ---------------------------------------------------------------------------
# 2. Add the new "metadata" sheet right after "data".
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Copy the header / index-column formatting from the "data" sheet so the new
# sheet's styled cells (bold, centered, bordered) reuse the same style.
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)     # xlPasteFormats
$excel.CutCopyMode = $false

# Header row.
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row.
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Palmoplantar Keratoderma and Erythrokeratoderma"
$metaSheet.Range("C2").Value = 153

# "0.107" must stay TEXT (not be auto-coerced to a number) but keep the
# default (unstyled) cell format, matching the source data. Temporarily force
# a text number-format so Excel doesn't coerce the input to a double, then
# clear the format override so the cell keeps the default style.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.107"
$metaSheet.Range("D2").ClearFormats()

$metaSheet.Range("E2").Value = "2021-09-15T09:54:26.416918Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:09.559036"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/153/?format=json"

# Restore "data" as the active sheet/selection (it was active before this
# script ran) so the only structural change is the new sheet itself.
$dataSheet.Activate()
